# Generate Report for Handback
#
# This script reflects a fresh localization-status report run:
#  - the zh-cn / de-de "Status" moves from "Ready for handoff" to
#    "Handed back: in sync with en-US"
#  - the "Latest Handback DateTime" timestamps are refreshed
#  - the stale "version mismatch" Error Detail message is cleared now that
#    the handback is in sync
#  - the Status / Error Detail columns are widened to fit the new text

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsOverview = $wb.Worksheets.Item("Overview")

# --- Status text: flips from "Ready for handoff" to the handed-back message
# on every sheet that surfaces it (per-language Status column plus the
# rolled-up Overview columns).
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# --- Latest Handback DateTime refreshed to the new generation run
$wsZhCn.Range("K2").Value = "2016-09-02 06:57:33"
$wsDeDe.Range("K2").Value = "2016-09-02 06:57:40"

# --- Error Detail no longer applicable now that handback is in sync
$wsZhCn.Range("P2").Value = ""
$wsDeDe.Range("P2").Value = ""

# --- Widen the Status / Error Detail columns to fit the longer text.
# Excel snaps ColumnWidth to whole-pixel (Calibri-11 MDW) increments, so the
# inputs below are chosen to land on the closest attainable grid width.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334
